$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (16 and 17) after the existing data (rows 2-15)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "2024-06-15 05:12:58"
$ws.Range("D16").Value = 200
$ws.Range("E16").Value = 8

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "2024-06-15 05:12:59"
$ws.Range("D17").Value = 200
$ws.Range("E17").Value = 2
